$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table used to end at row 75, which carried a bold "closing" bottom
# border (and taller row height). We are appending a new record (register
# 74 / 0x4A, "Silo full signal") below it, and the table no longer gets a
# special closing border at all - row 75 becomes a normal middle row and
# the new row 76 is also styled like a normal row.

# 1) Copy row 74's formatting (a normal, un-bolded middle row) down onto
#    both row 75 (to strip its old thick bottom border) and the new row 76.
$ws.Range("A74:C74").Copy()
$ws.Range("A75:C76").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 75 used to be taller (ht=15) to fit its bold closing border; now that
# it is a plain row again, let it fall back to the sheet's default height.
$ws.Rows.Item(75).AutoFit()
$ws.Rows.Item(76).AutoFit()

# 2) Populate the new row 76 with the new register.
$ws.Cells.Item(76, 1).Value2 = 74
$ws.Cells.Item(76, 2).Formula = "=DEC2HEX(A76,2)"
$ws.Cells.Item(76, 3).Value2 = "Silo full signal"

# 3) Update the selection to reflect where the author left the cursor.
$ws.Range("C79").Select()
